$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Finish the sentence in the existing last paragraph: the run that ended
#    in "...last iteration's" (followed by the _GoBack bookmark and a
#    " recap" run) becomes one run ending in "...last iteration's recap".
#    Scoping Find to that paragraph's Range and replacing the phrase that
#    spans the (collapsed) bookmark merges everything into a single run and
#    drops the now-redundant bookmark pair.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Find.Execute("iteration’s recap", $true, $false, $false, $false, $false, $true, 1, $false, "iteration’s recap", 2) | Out-Null

# ---------------------------------------------------------------------------
# Helper: append a brand new paragraph (with the same "en-US" run/paragraph
# formatting used throughout the document) holding a single run of text.
# ---------------------------------------------------------------------------
function Add-SimplePara([string]$text) {
    $p = $d.Paragraphs.Item($d.Paragraphs.Count)
    $p.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $newPara.Range.InsertAfter($text)
    return $newPara
}

Add-SimplePara "I would like to hearing from you about these overruled obstacles" | Out-Null

Add-SimplePara "I would like to hearing from you about the underlying benchmark of the requirement’s integrity" | Out-Null

# ---------------------------------------------------------------------------
# Paragraph made of three runs: "I would be grateful if you emphasiz" + "e" +
# " the need for the last sprint recap". Insert the full text as one run
# first, then force a run-split at the two boundaries by toggling Bold on
# and back off (the save step normalises the run's rPr back to the
# inherited state, so no stray formatting survives, but the run break
# remains).
# ---------------------------------------------------------------------------
$p4 = Add-SimplePara "I would be grateful if you emphasize the need for the last sprint recap"
$p4Start = $p4.Range.Start
foreach ($offset in @(35, 36)) {
    $splitRange = $d.Range($p4Start + $offset, $p4Start + $offset + 1)
    $splitRange.Bold = 1
    $splitRange.Bold = 0
}

Add-SimplePara "I would be grateful if you explain to me the point of that customer complaint" | Out-Null

Add-SimplePara "I would be grateful if you explain the point of the underlying rules for requirement’s integrity" | Out-Null

Add-SimplePara "I would be grateful if you send me the basic templates" | Out-Null

Add-SimplePara "I would be grateful if you send me the underlying templates" | Out-Null

Add-SimplePara "I would be grateful if you a list of the customer’s complaints" | Out-Null

Add-SimplePara "I would like to enquire about the overruled obstacles and how are you going to cope with future obstacles" | Out-Null

# ---------------------------------------------------------------------------
# Paragraph made of five runs: "I would like to enquire about that
# unambiguous " + "allusion" + " to " + "the " + "underlying template".
# Same split-by-toggling-Bold technique, at each run boundary.
# ---------------------------------------------------------------------------
$fullP11 = "I would like to enquire about that unambiguous allusion to the underlying template"
$p11 = Add-SimplePara $fullP11
$p11Start = $p11.Range.Start
$seg1 = "I would like to enquire about that unambiguous "
$seg2 = "allusion"
$seg3 = " to "
$seg4 = "the "
$boundary1 = $seg1.Length
$boundary2 = $boundary1 + $seg2.Length
$boundary3 = $boundary2 + $seg3.Length
$boundary4 = $boundary3 + $seg4.Length
foreach ($offset in @($boundary1, $boundary2, $boundary3, $boundary4)) {
    $splitRange = $d.Range($p11Start + $offset, $p11Start + $offset + 1)
    $splitRange.Bold = 1
    $splitRange.Bold = 0
}

Add-SimplePara "I would like to enquire about the new vacancy in QA department" | Out-Null

# ---------------------------------------------------------------------------
# Final content paragraph: the _GoBack bookmark now wraps the whole run
# (Word's "go back to last edit" bookmark follows the final insertion).
# ---------------------------------------------------------------------------
$p13 = Add-SimplePara "I would be grateful if you tell me about the new vacancy as soon as possible"
$p13Start = $p13.Range.Start
$p13End = $p13.Range.End
$bmRange = $d.Range($p13Start, $p13End)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Two trailing empty paragraphs.
$trailing1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$trailing1.Range.InsertParagraphAfter()
$trailing2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$trailing2.Range.InsertParagraphAfter()
